## AB - First 4 locations of row 4 added
##
## 1) Bump the "Update automatically" date placeholder cached text from
##    11/10/2016 -> 20/10/2016 everywhere it appears (slide master, every
##    slide layout, and the notes master).
## 2) Widen/retitle the two "Flood Plains" title boxes on slides 37 & 38
##    ("Flood Plains 2" / "Flood Plains 1").

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq "11/10/2016") {
            $shp.TextFrame.TextRange.Text = "20/10/2016"
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every slide layout off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($li)
}

# Notes master (its placeholder text doesn't persist through the
# TextFrame.TextRange path in this host, but the HeadersFooters
# accessor does take effect, so use that here instead)
$nm = $p.NotesMaster
if ($nm.Shapes.Item(2).TextFrame.TextRange.Text -eq "11/10/2016") {
    $nm.HeadersFooters.DateAndTime.Text = "20/10/2016"
}

function Find-ShapeByText($slide, $text) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq $text) {
            return $shp
        }
    }
    return $null
}

# --- Slide 37: "Flood Plains" -> "Flood Plains 2", box widened ---
$s37 = $p.Slides.Item(37)
$shp37 = Find-ShapeByText $s37 "Flood Plains"
$shp37.Left = 309.7498
$shp37.Width = 456.3302
$shp37.TextFrame.TextRange.Text = "Flood Plains 2"

# --- Slide 38: "Flood Plains" -> "Flood Plains 1", box widened ---
$s38 = $p.Slides.Item(38)
$shp38 = Find-ShapeByText $s38 "Flood Plains"
$shp38.Left = 309.7498
$shp38.Width = 450.5702
$shp38.TextFrame.TextRange.Text = "Flood Plains 1"
